$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "First trial ..." labels to "First-trial ..." (hyphenated)
$ws.Range("A13").Value = "First-trial volume "
$ws.Range("A14").Value = "First-trial duration"

# Move the active selection to A15
$ws.Range("A15").Select()
